{"js": "// Updates the worksheet date title and all 25 division-problem answers in the table,\n// in document order, matching each run's existing (old) text to find its target\n// before replacing it with the new text. Using getRange().insertText(..., replace)\n// on each paragraph/cell keeps the existing run/paragraph formatting intact.\n\nconst body = context.document.body;\n\n// 1) Title paragraph (first paragraph in the body).\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst titleOld = \"2023-08-30 Wednesday\";\nconst titleNew = \"2023-08-31 Thursday\";\nconst titlePara = paragraphs.items[0];\ntitlePara.load(\"text\");\nawait context.sync();\nif (titlePara.text.trim() === titleOld) {\n  titlePara.getRange().insertText(titleNew, Word.InsertLocation.replace);\n}\n\n// 2) Table of division answers: replace each cell's text by its (row, column)\n//    position, matched against the value already present there.\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\nconst cellEdits = [\n  { row: 0, col: 0, oldText: \"26\u00f78=3, 2\", newText: \"79\u00f72=39, 1\" },\n  { row: 0, col: 1, oldText: \"30\u00f73=10, 0\", newText: \"79\u00f78=9, 7\" },\n  { row: 0, col: 2, oldText: \"17\u00f76=2, 5\", newText: \"19\u00f76=3, 1\" },\n  { row: 0, col: 3, oldText: \"75\u00f75=15, 0\", newText: \"40\u00f78=5, 0\" },\n  { row: 0, col: 4, oldText: \"27\u00f76=4, 3\", newText: \"35\u00f79=3, 8\" },\n  { row: 4, col: 0, oldText: \"89\u00f79=9, 8\", newText: \"85\u00f72=42, 1\" },\n  { row: 4, col: 1, oldText: \"15\u00f73=5, 0\", newText: \"71\u00f75=14, 1\" },\n  { row: 4, col: 2, oldText: \"50\u00f74=12, 2\", newText: \"50\u00f79=5, 5\" },\n  { row: 4, col: 3, oldText: \"72\u00f78=9, 0\", newText: \"13\u00f74=3, 1\" },\n  { row: 4, col: 4, oldText: \"36\u00f76=6, 0\", newText: \"98\u00f72=49, 0\" },\n  { row: 8, col: 0, oldText: \"36\u00f79=4, 0\", newText: \"60\u00f75=12, 0\" },\n  { row: 8, col: 1, oldText: \"62\u00f73=20, 2\", newText: \"26\u00f75=5, 1\" },\n  { row: 8, col: 2, oldText: \"65\u00f78=8, 1\", newText: \"38\u00f74=9, 2\" },\n  { row: 8, col: 3, oldText: \"54\u00f78=6, 6\", newText: \"43\u00f74=10, 3\" },\n  { row: 8, col: 4, oldText: \"21\u00f78=2, 5\", newText: \"15\u00f73=5, 0\" },\n  { row: 12, col: 0, oldText: \"76\u00f75=15, 1\", newText: \"39\u00f73=13, 0\" },\n  { row: 12, col: 1, oldText: \"83\u00f74=20, 3\", newText: \"63\u00f79=7, 0\" },\n  { row: 12, col: 2, oldText: \"63\u00f76=10, 3\", newText: \"21\u00f78=2, 5\" },\n  { row: 12, col: 3, oldText: \"16\u00f78=2, 0\", newText: \"33\u00f79=3, 6\" },\n  { row: 12, col: 4, oldText: \"69\u00f75=13, 4\", newText: \"33\u00f75=6, 3\" },\n  { row: 16, col: 0, oldText: \"66\u00f75=13, 1\", newText: \"34\u00f79=3, 7\" },\n  { row: 16, col: 1, oldText: \"50\u00f72=25, 0\", newText: \"81\u00f79=9, 0\" },\n  { row: 16, col: 2, oldText: \"61\u00f77=8, 5\", newText: \"76\u00f75=15, 1\" },\n  { row: 16, col: 3, oldText: \"36\u00f73=12, 0\", newText: \"80\u00f77=11, 3\" },\n  { row: 16, col: 4, oldText: \"19\u00f76=3, 1\", newText: \"94\u00f78=11, 6\" },\n];\n\nconst cellRanges = cellEdits.map((edit) => {\n  const cell = table.getCell(edit.row, edit.col);\n  const range = cell.body.getRange();\n  range.load(\"text\");\n  return range;\n});\nawait context.sync();\n\ncellEdits.forEach((edit, i) => {\n  const range = cellRanges[i];\n  if (range.text.trim() === edit.oldText) {\n    range.insertText(edit.newText, Word.InsertLocation.replace);\n  }\n});\nawait context.sync();\n", "ps1": "# Update the worksheet date title and all 25 division-problem answers in the table.\n# Each cell/paragraph is located positionally (title paragraph, then table row/col)\n# and its current text is verified against the expected old value before the new\n# value is written, so formatting (fonts, alignment) on the existing run is kept.\n\n$d = $word.ActiveDocument\n\n# 1) Title paragraph (first paragraph in the document).\n$titleOld = \"2023-08-30 Wednesday\"\n$titleNew = \"2023-08-31 Thursday\"\n$titlePara = $d.Paragraphs.Item(1)\n$titleRange = $titlePara.Range\nif ($titleRange.Text.TrimEnd(\"`r\", \"`a\") -eq $titleOld) {\n    $titleRange.Text = $titleNew\n}\n\n# 2) Table of division answers: update each cell by its (row, column) position,\n#    matched against the value already present there.\n$table = $d.Tables.Item(1)\n\n$cellEdits = @(\n    @{ Row = 1; Col = 1; OldText = \"26\u00f78=3, 2\"; NewText = \"79\u00f72=39, 1\" }\n    @{ Row = 1; Col = 2; OldText = \"30\u00f73=10, 0\"; NewText = \"79\u00f78=9, 7\" }\n    @{ Row = 1; Col = 3; OldText = \"17\u00f76=2, 5\"; NewText = \"19\u00f76=3, 1\" }\n    @{ Row = 1; Col = 4; OldText = \"75\u00f75=15, 0\"; NewText = \"40\u00f78=5, 0\" }\n    @{ Row = 1; Col = 5; OldText = \"27\u00f76=4, 3\"; NewText = \"35\u00f79=3, 8\" }\n    @{ Row = 5; Col = 1; OldText = \"89\u00f79=9, 8\"; NewText = \"85\u00f72=42, 1\" }\n    @{ Row = 5; Col = 2; OldText = \"15\u00f73=5, 0\"; NewText = \"71\u00f75=14, 1\" }\n    @{ Row = 5; Col = 3; OldText = \"50\u00f74=12, 2\"; NewText = \"50\u00f79=5, 5\" }\n    @{ Row = 5; Col = 4; OldText = \"72\u00f78=9, 0\"; NewText = \"13\u00f74=3, 1\" }\n    @{ Row = 5; Col = 5; OldText = \"36\u00f76=6, 0\"; NewText = \"98\u00f72=49, 0\" }\n    @{ Row = 9; Col = 1; OldText = \"36\u00f79=4, 0\"; NewText = \"60\u00f75=12, 0\" }\n    @{ Row = 9; Col = 2; OldText = \"62\u00f73=20, 2\"; NewText = \"26\u00f75=5, 1\" }\n    @{ Row = 9; Col = 3; OldText = \"65\u00f78=8, 1\"; NewText = \"38\u00f74=9, 2\" }\n    @{ Row = 9; Col = 4; OldText = \"54\u00f78=6, 6\"; NewText = \"43\u00f74=10, 3\" }\n    @{ Row = 9; Col = 5; OldText = \"21\u00f78=2, 5\"; NewText = \"15\u00f73=5, 0\" }\n    @{ Row = 13; Col = 1; OldText = \"76\u00f75=15, 1\"; NewText = \"39\u00f73=13, 0\" }\n    @{ Row = 13; Col = 2; OldText = \"83\u00f74=20, 3\"; NewText = \"63\u00f79=7, 0\" }\n    @{ Row = 13; Col = 3; OldText = \"63\u00f76=10, 3\"; NewText = \"21\u00f78=2, 5\" }\n    @{ Row = 13; Col = 4; OldText = \"16\u00f78=2, 0\"; NewText = \"33\u00f79=3, 6\" }\n    @{ Row = 13; Col = 5; OldText = \"69\u00f75=13, 4\"; NewText = \"33\u00f75=6, 3\" }\n    @{ Row = 17; Col = 1; OldText = \"66\u00f75=13, 1\"; NewText = \"34\u00f79=3, 7\" }\n    @{ Row = 17; Col = 2; OldText = \"50\u00f72=25, 0\"; NewText = \"81\u00f79=9, 0\" }\n    @{ Row = 17; Col = 3; OldText = \"61\u00f77=8, 5\"; NewText = \"76\u00f75=15, 1\" }\n    @{ Row = 17; Col = 4; OldText = \"36\u00f73=12, 0\"; NewText = \"80\u00f77=11, 3\" }\n    @{ Row = 17; Col = 5; OldText = \"19\u00f76=3, 1\"; NewText = \"94\u00f78=11, 6\" }\n)\n\nforeach ($edit in $cellEdits) {\n    $cell = $table.Cell($edit.Row, $edit.Col)\n    $cellRange = $cell.Range\n    $currentText = $cellRange.Text.TrimEnd(\"`r\", \"`a\")\n    if ($currentText -eq $edit.OldText) {\n        $cellRange.Text = $edit.NewText\n    }\n}\n"}
